$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 92593120
$ws.Range("I92").Value = 4274111.5
$ws.Range("J92").Value = 666666700
$ws.Range("K92").Value = 4274111.5
$ws.Range("L92").Value = 666666700
$ws.Range("M92").Value = -4272863.5
$ws.Range("N92").Value = -666669196
$ws.Range("H96").Value = 1488.4
$ws.Range("I96").Value = 969
$ws.Range("J96").Value = 1834.6666
$ws.Range("K96").Value = 2907
$ws.Range("L96").Value = 5503.9998
$ws.Range("M96").Value = -1534
$ws.Range("N96").Value = -8249.9998
$ws.Range("H111").Value = 1810.909
$ws.Range("I111").Value = 1080.8889
$ws.Range("J111").Value = 5096
$ws.Range("K111").Value = 3242.6667
$ws.Range("L111").Value = 15288
$ws.Range("M111").Value = -175.6666999999998
$ws.Range("N111").Value = -21422
$ws.Range("H137").Value = 1682.5853
$ws.Range("I137").Value = 1363.1212
$ws.Range("J137").Value = 3000.375
$ws.Range("K137").Value = 4089.3636
$ws.Range("L137").Value = 9001.125
$ws.Range("M137").Value = -1539.3636
$ws.Range("N137").Value = -14101.125

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12208.171
$ws.Range("I32").Value = 8845.897999999999
$ws.Range("J32").Value = 20833.13
$ws.Range("K32").Value = 8845.897999999999
$ws.Range("L32").Value = 20833.13
$ws.Range("M32").Value = -8558.897999999999
$ws.Range("N32").Value = -21407.13
$ws.Range("H122").Value = 677736.6
$ws.Range("I122").Value = 857699.7
$ws.Range("J122").Value = 2875
$ws.Range("K122").Value = 2573099.1
$ws.Range("L122").Value = 8625
$ws.Range("M122").Value = -2570649.1
$ws.Range("N122").Value = -13525
$ws.Range("H132").Value = 1451686.6
$ws.Range("I132").Value = 1660.16
$ws.Range("J132").Value = 5267545.5
$ws.Range("K132").Value = 4980.48
$ws.Range("L132").Value = 15802636.5
$ws.Range("M132").Value = -2450.48
$ws.Range("N132").Value = -15807696.5
$ws.Range("H134").Value = 41121.5
$ws.Range("J134").Value = 41121.5
$ws.Range("L134").Value = 41121.5
$ws.Range("N134").Value = -51261.5
$ws.Range("H135").Value = 45745.668
$ws.Range("J135").Value = 45745.668
$ws.Range("L135").Value = 45745.668
$ws.Range("N135").Value = -55885.668
$ws.Range("H139").Value = 62647.5
$ws.Range("J139").Value = 62647.5
$ws.Range("L139").Value = 62647.5
$ws.Range("N139").Value = -72927.5

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 151848.2
$ws.Range("I107").Value = 189429
$ws.Range("J107").Value = 1525
$ws.Range("K107").Value = 189429
$ws.Range("L107").Value = 1525
$ws.Range("M107").Value = -187509
$ws.Range("N107").Value = -5365

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8890461
$ws.Range("I31").Value = 1589.4333
$ws.Range("J31").Value = 18414252
$ws.Range("K31").Value = 1589.4333
$ws.Range("L31").Value = 18414252
$ws.Range("M31").Value = -1294.4333
$ws.Range("N31").Value = -18414842
$ws.Range("H34").Value = 8890461
$ws.Range("I34").Value = 1589.4333
$ws.Range("J34").Value = 18414252
$ws.Range("K34").Value = 1589.4333
$ws.Range("L34").Value = 18414252
$ws.Range("M34").Value = -1387.4333
$ws.Range("N34").Value = -18414656
$ws.Range("H105").Value = 2214.6155
$ws.Range("I105").Value = 2280.9092
$ws.Range("J105").Value = 1850
$ws.Range("K105").Value = 2280.9092
$ws.Range("L105").Value = 1850
$ws.Range("M105").Value = -533.9092000000001
$ws.Range("N105").Value = -5344
$ws.Range("H110").Value = 42500
$ws.Range("J110").Value = 42500
$ws.Range("L110").Value = 42500
$ws.Range("N110").Value = -50680

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 4552475.5
$ws.Range("I5").Value = 11808
$ws.Range("J5").Value = 8267567.5
$ws.Range("K5").Value = 35424
$ws.Range("L5").Value = 24802702.5
$ws.Range("M5").Value = -35312
$ws.Range("N5").Value = -24802926.5
$ws.Range("H34").Value = 2670.5881
$ws.Range("J34").Value = 2985.7144
$ws.Range("L34").Value = 8957.143199999999
$ws.Range("N34").Value = -9125.143199999999
$ws.Range("H131").Value = 3572908.5
$ws.Range("I131").Value = 25000332
$ws.Range("J131").Value = 1671.0416
$ws.Range("K131").Value = 75000996
$ws.Range("L131").Value = 5013.1248
$ws.Range("M131").Value = -74995956
$ws.Range("N131").Value = -15093.1248
$ws.Range("H132").Value = 10057673
$ws.Range("I132").Value = 3002
$ws.Range("J132").Value = 11885795
$ws.Range("K132").Value = 27018
$ws.Range("L132").Value = 106972155
$ws.Range("M132").Value = -24488
$ws.Range("N132").Value = -106977215
$ws.Range("H135").Value = 4552475.5
$ws.Range("I135").Value = 11808
$ws.Range("J135").Value = 8267567.5
$ws.Range("K135").Value = 106272
$ws.Range("L135").Value = 74408107.5
$ws.Range("M135").Value = -103737
$ws.Range("N135").Value = -74413177.5

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 54594772
$ws.Range("I122").Value = 70989320
$ws.Range("J122").Value = 30002950
$ws.Range("K122").Value = 212967960
$ws.Range("L122").Value = 90008850
$ws.Range("M122").Value = -212965510
$ws.Range("N122").Value = -90013750
$ws.Range("H133").Value = 43525.715
$ws.Range("J133").Value = 43525.715
$ws.Range("L133").Value = 43525.715
$ws.Range("N133").Value = -53645.715

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 76925840
$ws.Range("I40").Value = 100002650
$ws.Range("J40").Value = 3151.6667
$ws.Range("K40").Value = 100002650
$ws.Range("L40").Value = 3151.6667
$ws.Range("M40").Value = -100002514
$ws.Range("N40").Value = -3423.6667
$ws.Range("H55").Value = 112.44444
$ws.Range("I55").Value = 123.666664
$ws.Range("K55").Value = 123.666664
$ws.Range("M55").Value = 49.333336
$ws.Range("H134").Value = 34954.11
$ws.Range("J134").Value = 34954.11
$ws.Range("L134").Value = 34954.11
$ws.Range("N134").Value = -45094.11
$ws.Range("H135").Value = 179751.6
$ws.Range("J135").Value = 179751.6
$ws.Range("L135").Value = 179751.6
$ws.Range("N135").Value = -189891.6
$ws.Range("H138").Value = 40354.668
$ws.Range("J138").Value = 40354.668
$ws.Range("L138").Value = 40354.668
$ws.Range("N138").Value = -50634.668

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 23906.572
$ws.Range("J86").Value = 23906.572
$ws.Range("L86").Value = 23906.572
$ws.Range("N86").Value = -26152.572
$ws.Range("H89").Value = 23906.572
$ws.Range("J89").Value = 23906.572
$ws.Range("L89").Value = 119532.86
$ws.Range("N89").Value = -130764.86
$ws.Range("H133").Value = 44747.75
$ws.Range("J133").Value = 44747.75
$ws.Range("L133").Value = 44747.75
$ws.Range("N133").Value = -54867.75
$ws.Range("H138").Value = 52535.6
$ws.Range("J138").Value = 54226
$ws.Range("L138").Value = 54226
$ws.Range("N138").Value = -64506
